$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Insert a new row of data into columns A:B only, at row 11 ---
# (Columns E,H,K,N,Q must NOT shift - only A/B shift down by one row.)
# Shift existing A11:B34 down to A12:B35, working from the bottom up.
for ($r = 34; $r -ge 11; $r--) {
    $srcA = $ws.Range("A$r").Text
    $srcB = $ws.Range("B$r").Text
    $dst = $r + 1
    if ($srcA -ne "") {
        $ws.Range("A$dst").Value = $srcA
    }
    if ($srcB -ne "") {
        $ws.Range("B$dst").Value = $srcB
    } else {
        $ws.Range("B$dst").ClearContents()
    }
}

# New row 11: a new question, with no similarity value in column B.
$ws.Range("A11").Value = "How does NAT (Network Address Translation) work?"
$ws.Range("B11").ClearContents()

# --- 2) Add new column T: "after prepro correlation ngram no len" ---
$tValues = @{}
$tValues[1] = "after prepro correlation ngram no len"
$tValues[2] = "['0.56', '0.52', '0.5']"
$tValues[3] = "['0.6', '0.39', '0.39']"
$tValues[4] = "['0.52', '0.41', '0.36']"
$tValues[5] = "['0.51', '0.37', '0.35']"
$tValues[6] = "['0.62', '0.57', '0.5']"
$tValues[7] = "['0.62', '0.5', '0.46']"
$tValues[8] = "['0.73', '0.55', '0.48']"
$tValues[9] = "['0.7', '0.5', '0.48']"
$tValues[10] = "['0.71', '0.45', '0.43']"
$tValues[11] = "['0.64', '0.46', '0.43']"
$tValues[12] = "['0.45', '0.39', '0.38']"
$tValues[13] = "['0.36', '0.35', '0.33']"
$tValues[14] = "['0.32', '0.29', '0.28']"
$tValues[15] = "['0.32', '0.27', '0.25']"
$tValues[16] = "['0.44', '0.42', '0.39']"
$tValues[17] = "['0.57', '0.49', '0.41']"
$tValues[18] = "['0.4', '0.39', '0.38']"
$tValues[19] = "['0.43', '0.41', '0.39']"
$tValues[20] = "['0.59', '0.56', '0.54']"
$tValues[21] = "['0.59', '0.5', '0.39']"
$tValues[22] = "['0.72', '0.65', '0.46']"
$tValues[23] = "['0.7', '0.41', '0.31']"
$tValues[24] = "['0.62', '0.54', '0.52']"
$tValues[25] = "['0.59', '0.45', '0.44']"
$tValues[26] = "['0.55', '0.44', '0.43']"
$tValues[27] = "['0.58', '0.54', '0.51']"
$tValues[28] = "['0.38', '0.35', '0.33']"
$tValues[29] = "['0.8', '0.47', '0.45']"
$tValues[30] = "['0.61', '0.43', '0.42']"
$tValues[31] = "['0.36', '0.29', '0.27']"
$tValues[32] = "['0.58', '0.44', '0.42']"
$tValues[33] = "['0.66', '0.47', '0.42']"
$tValues[34] = "['0.45', '0.43', '0.36']"
$tValues[35] = "['0.33', '0.29', '0.28']"

for ($r = 1; $r -le 35; $r++) {
    $ws.Range("T$r").Value = $tValues[$r]
}

# --- 3) Cosmetic view-state updates to mirror the edit session ---
$ws.Columns.Item(1).ColumnWidth = 52.6
$ws.Application.ActiveWindow.Zoom = 79
[void]$ws.Range("B11").Select()
